$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Clear the now-unused old "Nr." header cell in column C (the
#    table shifts one column to the right and gains a new column).
# ------------------------------------------------------------------
$ws.Range("C3").ClearContents()

# ------------------------------------------------------------------
# 2. Header row (row 3), columns D..H.
#    D3/E3/F3 already carry the bold "header" style (s=1); copy that
#    same formatting onto the new G3/H3 header cells before setting
#    their text.
# ------------------------------------------------------------------
$ws.Range("D3").Copy()
$ws.Range("G3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(3, 4).Value = "Nr."
$ws.Cells.Item(3, 5).Value = "Als…"
$ws.Cells.Item(3, 6).Value = "möchte ich …"
$ws.Cells.Item(3, 7).Value = "damit/weil/denn/um …"
$ws.Cells.Item(3, 8).Value = "Priorität "

# ------------------------------------------------------------------
# 3. Data rows 4..17, columns D..H.
# ------------------------------------------------------------------
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = "Benutzer"
$ws.Cells.Item(4, 6).Value = "Bewertungen anderer Benutzer lesen können"
$ws.Cells.Item(4, 7).Value = "sich entscheiden zu können, ob die Veranstaltung meinen Erwartungen entspricht"
$ws.Cells.Item(4, 8).Value = "S"
$ws.Cells.Item(5, 4).Value = 2
$ws.Cells.Item(5, 5).Value = "Benutzer"
$ws.Cells.Item(5, 6).Value = "mich Anmelden/Registrieren können"
$ws.Cells.Item(5, 7).Value = "personalisierte Funktionen nutzen zu können"
$ws.Cells.Item(5, 8).Value = "XL"
$ws.Cells.Item(6, 4).Value = 3
$ws.Cells.Item(6, 5).Value = "Benutzer"
$ws.Cells.Item(6, 6).Value = "meine Kontoinformationen bearbeiten können"
$ws.Cells.Item(6, 7).Value = "meine persönlichen Informationen zu aktualisieren oder zu korrigieren"
$ws.Cells.Item(6, 8).Value = "XL"
$ws.Cells.Item(7, 4).Value = 4
$ws.Cells.Item(7, 5).Value = "Eventorganisator"
$ws.Cells.Item(7, 6).Value = "Benachrichtigung an die Teilnehmer senden können"
$ws.Cells.Item(7, 7).Value = "auf wichtige aenderungen aufmerksam zu machen"
$ws.Cells.Item(7, 8).Value = "L"
$ws.Cells.Item(8, 4).Value = 5
$ws.Cells.Item(8, 5).Value = "Eventorganisator"
$ws.Cells.Item(8, 6).Value = "eine Liste der teilnehmenden Personen sehen können"
$ws.Cells.Item(8, 7).Value = "die Gästeliste zu verwalten und Ressourcen zu planen"
$ws.Cells.Item(8, 8).Value = "M"
$ws.Cells.Item(9, 4).Value = 6
$ws.Cells.Item(9, 5).Value = "Eventorganisator"
$ws.Cells.Item(9, 6).Value = "Veranstaltungsstandort auf einer Karte anzeigen"
$ws.Cells.Item(9, 7).Value = "eine visuelle Darstellung des Veranstaltungsorts zu haben"
$ws.Cells.Item(9, 8).Value = "S"
$ws.Cells.Item(10, 4).Value = 7
$ws.Cells.Item(10, 5).Value = "Eventorganisator"
$ws.Cells.Item(10, 6).Value = "eine Veranstaltung erstellen können"
$ws.Cells.Item(10, 7).Value = "alle relevanten Informationen zu erfassen und zu organisieren"
$ws.Cells.Item(10, 8).Value = "XL"
$ws.Cells.Item(11, 4).Value = 8
$ws.Cells.Item(11, 5).Value = "Eventorganisator"
$ws.Cells.Item(11, 6).Value = "Veranstaltungsdetails bearbeiten können"
$ws.Cells.Item(11, 7).Value = "aktualisierte Informationen an die Teilnehmer weiterzugeben"
$ws.Cells.Item(11, 8).Value = "XL"
$ws.Cells.Item(12, 4).Value = 9
$ws.Cells.Item(12, 5).Value = "Teilnehmer"
$ws.Cells.Item(12, 6).Value = "ein Event im Google Kalender speichern können "
$ws.Cells.Item(12, 7).Value = "den Veranstaltungstermin in meinem Kalender zu behalten"
$ws.Cells.Item(12, 8).Value = "M"
$ws.Cells.Item(13, 4).Value = 10
$ws.Cells.Item(13, 5).Value = "Teilnehmer"
$ws.Cells.Item(13, 6).Value = "einfach nach Events suchen können"
$ws.Cells.Item(13, 7).Value = "ich Veranstaltungen finde, die meinen Interessen entsprechen"
$ws.Cells.Item(13, 8).Value = "S"
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = "Teilnehmer"
$ws.Cells.Item(14, 6).Value = "alle Events durchsuchen können"
$ws.Cells.Item(14, 7).Value = "interessante Veranstaltungen zu finden"
$ws.Cells.Item(14, 8).Value = "S"
$ws.Cells.Item(15, 4).Value = 12
$ws.Cells.Item(15, 5).Value = "Teilnehmer"
$ws.Cells.Item(15, 6).Value = "mich zu events anmelden können"
$ws.Cells.Item(15, 7).Value = "an Veranstaltungen teilzunehmen"
$ws.Cells.Item(15, 8).Value = "XL"
$ws.Cells.Item(16, 4).Value = 13
$ws.Cells.Item(16, 5).Value = "Teilnehmer"
$ws.Cells.Item(16, 6).Value = "mich von angemeldeten Events abmelden können"
$ws.Cells.Item(16, 7).Value = "planänderungen an den Veranstalter mitzuteilen"
$ws.Cells.Item(16, 8).Value = "XL"
$ws.Cells.Item(17, 4).Value = 14
$ws.Cells.Item(17, 5).Value = "Teilnehmer"
$ws.Cells.Item(17, 6).Value = "Feedback zu einer vergangenen Veranstaltung geben"
$ws.Cells.Item(17, 7).Value = "die Qualität zukünftiger Events zu verbessern"
$ws.Cells.Item(17, 8).Value = "XS"

# ------------------------------------------------------------------
# 4. Grow the existing table ("Tabelle2") from C3:F11 to the new
#    D3:H17 range (5 columns x 14 data rows); the header-row text
#    already written above becomes each ListColumn's name.
# ------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("D3:H17"))

# ------------------------------------------------------------------
# 5. Cell-level formatting tweaks.
#    - Column D (the "Nr." numbers) is left-aligned.
#    - H8:H10 also picked up the left-aligned style in the authored
#      file (same visual style as column D, just also applied here).
#    - G16 lost its inherited theme color (now plain Calibri 11).
# ------------------------------------------------------------------
$ws.Range("D4:D17").HorizontalAlignment = -4131   # xlLeft
$ws.Range("H8:H10").HorizontalAlignment = -4131   # xlLeft

$ws.Range("G16").Font.Name = "Calibri"
$ws.Range("G16").Font.Size = 11

# ------------------------------------------------------------------
# 6. Column widths (characters) for the now-used columns C..H.
# ------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 7.6
$ws.Columns.Item(4).ColumnWidth = 8.3
$ws.Columns.Item(5).ColumnWidth = 23.6
$ws.Columns.Item(6).ColumnWidth = 63.75
$ws.Columns.Item(7).ColumnWidth = 86.9
$ws.Columns.Item(8).ColumnWidth = 12

# ------------------------------------------------------------------
# 7. Selection / scroll position, matching the saved view.
# ------------------------------------------------------------------
$ws.Range("G18").Select()
